$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.643.47'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '3.107.67'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '523.83'
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.11'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.107.71'
$ws.Range("E8").Value = '  +1.89%  '
$ws.Range("E9").Value = '  +0.30%  '
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("E11").Value = '  +2.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.386'
$ws.Range("E12").Value = '  +3.14%  '
$ws.Range("D13").Value = '3.640.95'
$ws.Range("E13").Value = '  +1.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.132'
$ws.Range("E14").Value = '  +1.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.25'
$ws.Range("E15").Value = '  +2.88%  '
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("D17").Value = '57.704.30'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = '3.107.20'
$ws.Range("E18").Value = '  +1.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.13'
$ws.Range("E19").Value = '  +1.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.85'
$ws.Range("E20").Value = '  +0.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.09'
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '336.49'
$ws.Range("E22").Value = '  +2.11%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.514'
$ws.Range("E24").Value = '  +3.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.56'
$ws.Range("E25").Value = '  +1.30%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("E28").Value = '  +3.15%  '
$ws.Range("E29").Value = '  +3.96%  '
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("E31").Value = '  +1.18%  '
$ws.Range("E32").Value = '  +2.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.02'
$ws.Range("E34").Value = '  +2.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '155.00'
$ws.Range("E35").Value = '  +0.33%  '
$ws.Range("E36").Value = '  +4.55%  '
$ws.Range("E37").Value = '  +3.47%  '
$ws.Range("E38").Value = '  -0.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.29'
$ws.Range("E39").Value = '  +2.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0665'
$ws.Range("E40").Value = '  -0.86%  '
$ws.Range("B41").Value = 'RenzoRestakedETH'
$ws.Range("C41").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D41").Value = '3.148.93'
$ws.Range("E41").Value = '  +1.87%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.687'
$ws.Range("E42").Value = '  +5.79%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.52'
$ws.Range("E43").Value = '  +12.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.94'
$ws.Range("E44").Value = '  +1.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '36.75'
$ws.Range("E45").Value = '  +0.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '2.314.47'
$ws.Range("E47").Value = '  +2.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0260'
$ws.Range("E48").Value = '  +1.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.978'
$ws.Range("E49").Value = '  +6.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.76'
$ws.Range("E50").Value = '  +0.38%  '
$ws.Range("E51").Value = '  +2.92%  '
